# "Saya edit backend tabel destinasi dan pengguna"
# Rework the "Data dokter" sheet into a "Data pengguna" (users) table:
#  - header row: ID / Nama Pengguna / Password / Email / Tanggal Lahir
#  - 3 data rows of user records (was 4 rows of doctor records)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1 stays "ID") --------------------------------------
$ws.Range("B1").Value = "Nama Pengguna"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "Tanggal Lahir"

# --- Remove the old 4th/5th data rows we don't need -------------------
# (only 3 user rows remain after the edit; drop the extra row first so
#  row numbers 2-4 are the ones we repopulate below)
$ws.Rows(5).Delete()

# --- Data rows ----------------------------------------------------------
# Column E holds dates stored as literal text (e.g. "2003-04-07"); format
# those cells as Text first so Excel doesn't silently convert the literal
# into a date serial number.

# Row 2 - Arslan Wiratama
$ws.Range("A2").Value = 12211845
$ws.Range("B2").Value = "Arslan Wiratama"
$ws.Range("C2").Value = "d41d8cd98f00b204e9800998ecf8427e"
$ws.Range("D2").Value = "izzuddinalfatah33@gmail.com"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2003-04-07"

# Row 3 - Yusril Camelia
$ws.Range("A3").Value = 12211846
$ws.Range("B3").Value = "Yusril Camelia"
$ws.Range("C3").Value = "d41d8cd98f00b204e9800998ecf8427e"
$ws.Range("D3").Value = "izzalfatah347@gmail.com"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2004-09-22"

# Row 4 - Ravin Sadewa
$ws.Range("A4").Value = 12211847
$ws.Range("B4").Value = "Ravin Sadewa"
$ws.Range("C4").Value = "d41d8cd98f00b204e9800998ecf8427e"
$ws.Range("D4").Value = "muhalfatah743@gmail.com"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2005-01-11"
